# Refresh the cryptos list values (GitHub Actions run, Wed Mar 13 14:54:57 UTC 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '72.459.00'
$ws.Range('E2').Value = '  +0.84%  '

# Row 3
$ws.Range('D3').Value = '3.971.04'
$ws.Range('E3').Value = '  -0.50%  '

# Row 4
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Range('E4').Value = '  +0.12%  '

# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '584.36'
$cell.ClearFormats()
$ws.Range('E5').Value = '  +7.89%  '

# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '153.73'
$cell.ClearFormats()
$ws.Range('E6').Value = '  +2.98%  '

# Row 7
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.674'
$cell.ClearFormats()
$ws.Range('E7').Value = '  -3.03%  '

# Row 8
$ws.Range('E8').Value = '  -0.08%  '

# Row 9
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.742'
$cell.ClearFormats()
$ws.Range('E9').Value = '  +0.04%  '

# Row 10
$ws.Range('E10').Value = '  -1.86%  '

# Row 11
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '52.75'
$cell.ClearFormats()
$ws.Range('E11').Value = '  +5.81%  '

# Row 12
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.0000316'
$cell.ClearFormats()
$ws.Range('E12').Value = '  -1.70%  '

# Row 13
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '10.70'
$cell.ClearFormats()
$ws.Range('E13').Value = '  +0.49%  '

# Row 14
$ws.Range('D14').Value = '4.605.97'
$ws.Range('E14').Value = '  -0.63%  '

# Row 15
$ws.Range('D15').Value = '3.977.08'
$ws.Range('E15').Value = '  -0.39%  '

# Row 16
$ws.Range('E16').Value = '  +7.60%  '

# Row 17
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '13.90'
$cell.ClearFormats()
$ws.Range('E17').Value = '  -0.86%  '

# Row 18
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '20.35'
$cell.ClearFormats()
$ws.Range('E18').Value = '  -0.04%  '

# Row 19
$ws.Range('E19').Value = '  -0.12%  '

# Row 20
$ws.Range('D20').Value = '72.345.60'
$ws.Range('E20').Value = '  +0.82%  '

# Row 21
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '425.67'
$cell.ClearFormats()
$ws.Range('E21').Value = '  -0.25%  '

# Row 22
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '4.65'
$cell.ClearFormats()
$ws.Range('E22').Value = '  +9.68%  '

# Row 23
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '95.16'
$cell.ClearFormats()
$ws.Range('E23').Value = '  -1.86%  '

# Row 24
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '3.42'
$cell.ClearFormats()
$ws.Range('E24').Value = '  -1.71%  '

# Row 25
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '14.15'
$cell.ClearFormats()
$ws.Range('E25').Value = '  -0.38%  '

# Row 26
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '4.38'
$cell.ClearFormats()
$ws.Range('E26').Value = '  +19.01%  '

# Row 27
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '11.18'
$cell.ClearFormats()
$ws.Range('E27').Value = '  -0.70%  '

# Row 28
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '10.66'
$cell.ClearFormats()
$ws.Range('E28').Value = '  +0.15%  '

# Row 29
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '5.92'
$cell.ClearFormats()
$ws.Range('E29').Value = '  +1.29%  '

# Row 30
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '36.12'
$cell.ClearFormats()
$ws.Range('E30').Value = '  -1.56%  '

# Row 31
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '7.67'
$cell.ClearFormats()
$ws.Range('E31').Value = '  +6.22%  '

# Row 32
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '13.36'
$cell.ClearFormats()
$ws.Range('E32').Value = '  +0.40%  '

# Row 33
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '49.11'
$cell.ClearFormats()
$ws.Range('E33').Value = '  +1.38%  '

# Row 34
$ws.Range('E34').Value = '  -0.73%  '

# Row 35
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '676.50'
$cell.ClearFormats()
$ws.Range('E35').Value = '  +0.51%  '

# Row 36
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '67.36'
$cell.ClearFormats()
$ws.Range('E36').Value = '  +3.30%  '

# Row 37
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.435'
$cell.ClearFormats()
$ws.Range('E37').Value = '  -1.26%  '

# Row 38
$ws.Range('D38').Value = '0.0₃0849'
$ws.Range('E38').Value = '  +4.05%  '

# Row 39
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '3.35'
$cell.ClearFormats()
$ws.Range('E39').Value = '  +0.33%  '

# Row 40
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Range('E40').Value = '  +0.02%  '

# Row 41
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '0.145'
$cell.ClearFormats()
$ws.Range('E41').Value = '  -3.18%  '

# Row 42
$ws.Range('E42').Value = '  +0.20%  '

# Row 43
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '3.28'
$cell.ClearFormats()
$ws.Range('E43').Value = '  -3.10%  '

# Row 44
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '10.86'
$cell.ClearFormats()
$ws.Range('E44').Value = '  +10.78%  '

# Row 45
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.0482'
$cell.ClearFormats()
$ws.Range('E45').Value = '  -0.75%  '

# Row 46
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '2.71'
$cell.ClearFormats()
$ws.Range('E46').Value = '  +2.44%  '

# Row 47
$ws.Range('E47').Value = '  -1.24%  '

# Row 48
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '3.36'
$cell.ClearFormats()
$ws.Range('E48').Value = '  +0.41%  '

# Row 49
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '3.40'
$cell.ClearFormats()
$ws.Range('E49').Value = '  +5.11%  '

# Row 50
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '2.14'
$cell.ClearFormats()
$ws.Range('E50').Value = '  +7.28%  '

# Row 51
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '2.97'
$cell.ClearFormats()
$ws.Range('E51').Value = '  -0.93%  '
